# Apply the StructureDefinition-reference-path update:
#  - Version bump 5.0.0 -> 6.0.0
#  - Date bump
#  - Publisher filled in ("Alvearie Team")
#  - Contact/"No display for ContactDetail" row replaced by Jurisdiction/"United States of America"
#  - The now-redundant duplicate "Contact" row removed entirely
#  - Elements sheet: Extension's Short/Definition text corrected

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "6.0.0"
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$meta.Range("B9").Value = "Alvearie Team"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Old row 11 was a duplicate "Contact / No display for ContactDetail" row;
# deleting it shifts every subsequent row up by one, matching the new layout.
$meta.Rows.Item(11).Delete()

$elements = $wb.Worksheets.Item("Elements")

$elements.Range("K2").Value = "Reference Path"
$elements.Range("L2").Value = "Path to FHIR element in the reference that was leveraged to produce the insight."
